$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999434228881,
    0.99826174156708991,
    0.99515616211253743,
    0.99872574874011355,
    0.98916435770263056,
    0.96813710547393983,
    0.96412200137206749,
    0.95890604309932181,
    0.95386990388721715,
    0.94975644189193198,
    0.9481606427293916,
    0.94538137079262863,
    0.93409613795980462,
    0.92992916958025928,
    0.92733786865454937,
    0.92483150607213094,
    0.92112370011351297,
    0.92001478923547142,
    0.99645163170266027,
    0.98235065719314396,
    0.97926605663889132,
    0.97530051041104049,
    0.98812209487843372,
    0.97510196481348377,
    0.96864509250660191,
    0.96570809874385533,
    0.96257092011383816,
    0.95184758322152174,
    0.94454727350792456,
    0.94077262375797122,
    0.93312004400819049,
    0.93144075884011701,
    0.9309207475010175
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
